$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  46"
$ws.Range("C9").Value = "Report Covering the Week  11/13/2023  Through  11/19/2023"

# --- Plain value changes (style unchanged) ---
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 13
$ws.Range("K15").Value = 8.333333333333
$ws.Range("L15").Value = -23.529411764705
$ws.Range("M15").Value = -59.375
$ws.Range("N15").Value = -74.509803921568
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 50
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 18.181818181818
$ws.Range("I16").Value = 141
$ws.Range("J16").Value = 99
$ws.Range("K16").Value = 42.424242424242
$ws.Range("L16").Value = 80.769230769230
$ws.Range("M16").Value = -45.559845559845
$ws.Range("N16").Value = -84.471365638766
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 8
$ws.Range("E17").Value = -37.5
$ws.Range("F17").Value = 33
$ws.Range("G17").Value = 32
$ws.Range("H17").Value = 3.125
$ws.Range("I17").Value = 414
$ws.Range("J17").Value = 356
$ws.Range("K17").Value = 16.292134831460
$ws.Range("L17").Value = 71.074380165289
$ws.Range("M17").Value = 51.648351648351
$ws.Range("N17").Value = -41.359773371104
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 12
$ws.Range("H18").Value = 9.090909090909
$ws.Range("I18").Value = 116
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = 14.851485148514
$ws.Range("L18").Value = 114.814814814815
$ws.Range("M18").Value = -57.818181818181
$ws.Range("N18").Value = -92.472420506164
$ws.Range("C19").Value = 6
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 50
$ws.Range("F19").Value = 32
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = 14.285714285714
$ws.Range("I19").Value = 375
$ws.Range("J19").Value = 305
$ws.Range("K19").Value = 22.950819672131
$ws.Range("L19").Value = 53.688524590163
$ws.Range("M19").Value = 1.902173913043
$ws.Range("N19").Value = -28.435114503816
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 50
$ws.Range("G20").Value = 17
$ws.Range("H20").Value = -29.411764705882
$ws.Range("I20").Value = 141
$ws.Range("J20").Value = 111
$ws.Range("K20").Value = 27.027027027027
$ws.Range("L20").Value = 113.636363636364
$ws.Range("M20").Value = -20.338983050847
$ws.Range("N20").Value = -88.279301745635
$ws.Range("C21").Value = 26
$ws.Range("D21").Value = 22
$ws.Range("E21").Value = 18.181818181818
$ws.Range("F21").Value = 103
$ws.Range("H21").Value = 3
$ws.Range("I21").Value = 1212
$ws.Range("J21").Value = 990
$ws.Range("K21").Value = 22.424242424242
$ws.Range("L21").Value = 71.186440677966
$ws.Range("M21").Value = -13.118279569892
$ws.Range("N21").Value = -75.515151515151
$ws.Range("F23").Value = 4
$ws.Range("G23").Value = 2
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 74
$ws.Range("K23").Value = 48
$ws.Range("L23").Value = 60.869565217391
$ws.Range("M23").Value = 76.190476190476
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 10
$ws.Range("E24").Value = 260
$ws.Range("F24").Value = 131
$ws.Range("G24").Value = 84
$ws.Range("H24").Value = 55.952380952380
$ws.Range("I24").Value = 1141
$ws.Range("J24").Value = 1049
$ws.Range("K24").Value = 8.770257387988
$ws.Range("L24").Value = 56.515775034293
$ws.Range("M24").Value = -15.481481481481
$ws.Range("C25").Value = 16
$ws.Range("D25").Value = 17
$ws.Range("E25").Value = -5.882352941176
$ws.Range("F25").Value = 57
$ws.Range("G25").Value = 40
$ws.Range("H25").Value = 42.5
$ws.Range("I25").Value = 651
$ws.Range("J25").Value = 563
$ws.Range("K25").Value = 15.630550621669
$ws.Range("L25").Value = 32.317073170731
$ws.Range("M25").Value = -43.193717277486
$ws.Range("F26").Value = 3
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 22
$ws.Range("K26").Value = -26.666666666666
$ws.Range("L26").Value = -31.25
$ws.Range("L27").Value = 66
$ws.Range("J28").Value = 25
$ws.Range("K28").Value = -12
$ws.Range("L28").Value = -18.518518518518
$ws.Range("N28").Value = -78.217821782178
$ws.Range("J29").Value = 18
$ws.Range("K29").Value = 11.111111111111
$ws.Range("L29").Value = -16.666666666666
$ws.Range("N29").Value = -77.528089887640

# --- Style + value changes (need format paste from donor cells) ---
$ws.Range("C15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D15").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E15").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("F15").Value = 1
$ws.Range("I14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D23").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E23").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E23").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C26").Value = 2
$ws.Range("I14").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D26").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E26").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("C27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D27").Value = "'0"
$ws.Range("A14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E27").Value = "'***.*"
$ws.Range("A14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E28").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G28").Value = 2
$ws.Range("I14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H28").Value = 0
$ws.Range("K14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("D29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("D29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("E29").Value = -100
$ws.Range("K14").Copy()
$ws.Range("E29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G29").Value = 1
$ws.Range("I14").Copy()
$ws.Range("G29").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("H29").Value = 100
$ws.Range("K14").Copy()
$ws.Range("H29").PasteSpecial(-4122)
$excel.CutCopyMode = 0